$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 44 — this shifts the existing rows 44..129
# down to 45..130 (matching the target dimension A1:R130).
$ws.Rows(44).Insert()

# Populate the newly inserted row 44 with its data.
$ws.Cells.Item(44, 1).Value = 7
$ws.Cells.Item(44, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(44, 3).Value = "Ñuble"
$ws.Cells.Item(44, 4).Value = 44477
$ws.Cells.Item(44, 5).Value = 16
$ws.Cells.Item(44, 6).Value = 100112006
$ws.Cells.Item(44, 7).Value = "Repollo"
$ws.Cells.Item(44, 8).Value = "Crespo record"
$ws.Cells.Item(44, 9).Value = "Primera"
$ws.Cells.Item(44, 10).Value = 300
$ws.Cells.Item(44, 11).Value = 600
$ws.Cells.Item(44, 12).Value = 650
$ws.Cells.Item(44, 13).Value = 625
$ws.Cells.Item(44, 14).Value = "$/unidad"
$ws.Cells.Item(44, 15).Value = "Provincia de Diguillín"
$ws.Cells.Item(44, 16).Value = 625
$ws.Cells.Item(44, 17).Value = 1
$ws.Cells.Item(44, 18).Value = "Hortaliza"
